$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

function Set-CellXml($row, $col, $innerXml) {
    $cell = $t.Cell($row, $col)
    $cell.Range.InsertXML("<w:p $wNs>$innerXml</w:p>")
}

# --- Row 3, Col 4: collapse the spell-checked "textfield"/"passwordfield"
#     runs (w:proofErr wrappers) into a single plain run ---
Set-CellXml 3 4 '<w:r><w:t>Replaced textfield with a passwordfield</w:t></w:r>'

# --- Row 8, Col 4: drop the _GoBack bookmark that used to sit here,
#     text itself is unchanged ---
Set-CellXml 8 4 '<w:r><w:t>User logins, logs out, registers or them accessing content panels gets logged into a .txt file.</w:t></w:r>'

# --- Row 9: previously four empty cells ---
Set-CellXml 9 1 '<w:r><w:t>Passwords are not encrypted</w:t></w:r>'
Set-CellXml 9 2 '<w:r><w:t>Main.java, SQLite.java</w:t></w:r>'
Set-CellXml 9 3 '<w:r><w:t>Passwords can be easily hacked</w:t></w:r>'
Set-CellXml 9 4 '<w:r><w:t xml:space="preserve">Use </w:t></w:r><w:r><w:t>pbkdf2withhmacsha1</w:t></w:r><w:r><w:t xml:space="preserve"> to hash passwords with salt</w:t></w:r>'

# --- Row 10: previously four empty cells ---
Set-CellXml 10 1 '<w:r><w:t>Newly registered users get deleted when program restarts</w:t></w:r>'
Set-CellXml 10 2 '<w:r><w:t>Main.java</w:t></w:r>'
Set-CellXml 10 3 '<w:r><w:t>Integrity of data is compromised</w:t></w:r>'
Set-CellXml 10 4 '<w:r><w:t>Do not dropUserTable every time the program initializes</w:t></w:r>'

# --- Row 11: previously four empty cells; the _GoBack bookmark now
#     lives at the end of the last cell ---
Set-CellXml 11 1 '<w:r><w:t>People can guess account details an infinite amount of times</w:t></w:r>'
Set-CellXml 11 2 '<w:r><w:t>SQLite.java, LogWrite.java</w:t></w:r>'
Set-CellXml 11 3 '<w:r><w:t>Account details can be hacked/guessed through brute force</w:t></w:r>'
Set-CellXml 11 4 '<w:r><w:t>Put an attempt counter and disable logins when a certain number is reached</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>'
